$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# Row 2 (the "15/09/2005" row) -> date 10/12/2015, version 0.1, plus the
# previously-empty "Mo ta thay doi", "Nguoi lap" and "Nguoi duyet" cells.
# ---------------------------------------------------------------------------

# Col 1: Ngay lap
$d.Content.Find.Execute("15/09/2005", $true, $false, $false, $false, $false, $true, 1, $false, "10/12/2015", 2) | Out-Null

# Col 2: Mo ta thay doi (two new paragraphs)
$cell = $t.Cell(2, 2).Range
$cell.Collapse(0)
$cell.InsertBefore("Hoàn thành các tính năng cơ bản.`rPhiên bản Beta")

# Col 3: Phien ban  0.4 -> 0.1
$d.Content.Find.Execute("0.4", $true, $false, $false, $false, $false, $true, 1, $false, "0.1", 2) | Out-Null

# Col 4: Nguoi lap
$cell = $t.Cell(2, 4).Range
$cell.Collapse(0)
$cell.InsertBefore("Nguyễn Tiến Phong")

# Col 5: Nguoi duyet
$cell = $t.Cell(2, 5).Range
$cell.Collapse(0)
$cell.InsertBefore("Lê Anh Dũng")

# ---------------------------------------------------------------------------
# Row 3 (the "15/11/2005" row) -> date 15/12/2015, version 0.1.1, plus the
# previously-empty "Mo ta thay doi", "Nguoi lap" and "Nguoi duyet" cells.
# ---------------------------------------------------------------------------

# Col 1: Ngay lap
$d.Content.Find.Execute("15/11/2005", $true, $false, $false, $false, $false, $true, 1, $false, "15/12/2015", 2) | Out-Null

# Col 2: Mo ta thay doi (three new paragraphs)
$cell = $t.Cell(3, 2).Range
$cell.Collapse(0)
$text = "- Khắc phục lỗi khi nhiều người truy cập ứng dụng cùng một lúc" + "`r" + `
         "- Khắc phục lỗi Crash To Desktop (CTD) ngẫu nhiên khi qua được 1 màn chơi" + "`r" + `
         "- Khắc phục lỗi CTD khi người dùng vừa thoát ứng dụng"
$cell.InsertBefore($text)

# Col 3: Phien ban  0.6 -> 0.1.1
$d.Content.Find.Execute("0.6", $true, $false, $false, $false, $false, $true, 1, $false, "0.1.1", 2) | Out-Null

# Col 4: Nguoi lap
$cell = $t.Cell(3, 4).Range
$cell.Collapse(0)
$cell.InsertBefore("Nguyễn Thành Long")

# Col 5: Nguoi duyet
$cell = $t.Cell(3, 5).Range
$cell.Collapse(0)
$cell.InsertBefore("Lê Anh Dũng")
